$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cách cục (pattern) strings used in the added rows.
$s126 = "Sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Kình Dương, Đà La"
$s127 = "Sinh năm Đinh có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Kình Dương, Đà La"
$s128 = "Sinh năm Kỷ có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Kình Dương, Đà La"
$s129 = "Nam mệnh sinh năm Nhâm có Tử Vi tọa thủ cung Mệnh ở Hợi"
$s130 = "Nam mệnh sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Hợi"

# Row 69
$ws.Range("A69").Value = $s126
$ws.Range("B69").Value = $s126

# Row 70
$ws.Range("A70").Value = $s126
$ws.Range("B70").Value = $s127

# Row 71
$ws.Range("A71").Value = $s128
$ws.Range("B71").Value = $s128

# Row 72
$ws.Range("A72").Value = $s129
$ws.Range("B72").Value = $s129

# Row 73
$ws.Range("A73").Value = $s130
$ws.Range("B73").Value = $s130

# Match the existing "highlight" style (yellow fill) used across column A
# (and column B for rows where both cells repeat the same text).
$ws.Range("A69:B73").Interior.Color = 65535

# Update selection to match the post-edit state captured in the workbook.
$ws.Range("B72").Select()
